$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").Value = "2025/12/06 01:00"
$ws.Range("B79").Value = "-"
$ws.Range("C79").Value = "-"
$ws.Range("D79").Value = "-"
$ws.Range("E79").Value = "-"
$ws.Range("F79").Value = "-"
$ws.Range("G79").Value = "-"
